$d = $word.ActiveDocument

# --- Edit 1: merge "Kapil " / "Suryawanshi" / " (2039409)" runs into one run ---
$d.Content.Find.Execute("Kapil Suryawanshi (2039409)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Kapil Suryawanshi (2039409)", 2)

# --- Edit 2: add a new paragraph after the "Target users include..." paragraph ---
$r = $d.Content
$r.Find.Execute("Target users include individuals, families, students, and shared households seeking effective financial management tools. The usability goals of the " + [char]34 + "Expense Tracker" + [char]34 + " focus on providing an easy-to-use platform that offers transparent views of finances, personalized experiences tailored to individual needs, insightful visual reports, and collaborative financial management features.", `
                 $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter([char]13 + "No changes have been made to the project specification and design since the start of the project. The priority goals, features and technology stack are carefully tracked and applied throughout the development lifecycle.")

# Keep the document's TrackRevisions setting as authored, but make sure none of
# the edits above linger as visible tracked-change markup in the saved content.
$d.AcceptAllRevisions()
